# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# to the latest scraped values. Price cells use a leading apostrophe so Excel
# keeps them as literal text (matching the source data, which stores prices
# as plain strings such as "62.866.46" or "0.0902") instead of silently
# reinterpreting them as numbers, which would both drop the original
# formatting (e.g. thousand-separator dots) and introduce floating point
# rounding artifacts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.866.46'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '''2.466.30'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''571.84'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '''148.05'
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''0.530'
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("E10").Value = '  -0.38%  '
$ws.Range("D11").Value = '''5.19'
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").Value = '''0.348'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").Value = '''28.84'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("E14").Value = '  -2.62%  '
$ws.Range("D15").Value = '''2.918.84'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '''62.792.14'
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").Value = '''2.469.57'
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("D18").Value = '''7.66'
$ws.Range("E18").Value = '  -6.87%  '
$ws.Range("E19").Value = '  -2.85%  '
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").Value = '''322.12'
$ws.Range("E21").Value = '  -2.35%  '
$ws.Range("D22").Value = '''4.15'
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '''10.17'
$ws.Range("E24").Value = '  +3.52%  '
$ws.Range("D25").Value = '''64.89'
$ws.Range("E25").Value = '  -2.05%  '
$ws.Range("D26").Value = '''645.16'
$ws.Range("E26").Value = '  -4.15%  '
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("D28").Value = '''0.0₃0967'
$ws.Range("E28").Value = '  -3.25%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("E30").Value = '  -3.14%  '
$ws.Range("D31").Value = '''7.89'
$ws.Range("E31").Value = '  -2.51%  '
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  -3.79%  '
$ws.Range("E36").Value = '  -2.93%  '
$ws.Range("D37").Value = '''5.38'
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = '''0.365'
$ws.Range("E38").Value = '  -2.05%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '''150.13'
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("E42").Value = '  -2.18%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  -2.56%  '
$ws.Range("D45").Value = '''152.98'
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("D48").Value = '''20.26'
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").Value = '''0.605'
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D51").Value = '''0.0902'
$ws.Range("E51").Value = '  -1.84%  '
